# Update cryptocurrency Price (D) and Volume(1h) (E) columns
# with refreshed values from the data source (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'288.77"
$ws.Range("E2").Value = "'0.16%"
$ws.Range("D3").Value = "'31.03"
$ws.Range("D4").Value = "'4.959"
$ws.Range("E4").Value = "'0.47%"
$ws.Range("D5").Value = "'0.07355"
$ws.Range("E5").Value = "'1.59%"
$ws.Range("D6").Value = "'2.354"
$ws.Range("E6").Value = "'32.34%"
$ws.Range("D7").Value = "'7.731"
$ws.Range("E7").Value = "'1.83%"
$ws.Range("D8").Value = "'0.9124"
$ws.Range("E8").Value = "'1.20%"
$ws.Range("D9").Value = "'0.09380"
$ws.Range("E9").Value = "'20.45%"
$ws.Range("D10").Value = "'0.1703"
$ws.Range("E10").Value = "'2.06%"
$ws.Range("D11").Value = "'0.08294"
$ws.Range("E11").Value = "'4.38%"
$ws.Range("D12").Value = "'0.03113"
$ws.Range("E12").Value = "'1.79%"
$ws.Range("D13").Value = "'0.09967"
$ws.Range("E13").Value = "'-0.53%"
$ws.Range("D14").Value = "'0.001494"
$ws.Range("E14").Value = "'-0.33%"
$ws.Range("D15").Value = "'0.005779"
$ws.Range("E15").Value = "'0.53%"
$ws.Range("D16").Value = "'3.469"
$ws.Range("E16").Value = "'-0.12%"
$ws.Range("D17").Value = "'3.739"
$ws.Range("E17").Value = "'0.67%"
$ws.Range("D18").Value = "'2.008"
$ws.Range("E18").Value = "'-3.30%"
$ws.Range("D19").Value = "'0.3322"
$ws.Range("E19").Value = "'0.17%"
$ws.Range("D20").Value = "'0.1288"
$ws.Range("E20").Value = "'-0.56%"
$ws.Range("D21").Value = "'4.150"
$ws.Range("E21").Value = "'4.38%"
$ws.Range("E22").Value = "'-3.68%"
$ws.Range("D23").Value = "'0.04515"
$ws.Range("E23").Value = "'0.06%"
$ws.Range("D24").Value = "'0.001214"
$ws.Range("D25").Value = "'0.004178"
$ws.Range("E25").Value = "'-8.54%"
$ws.Range("D26").Value = "'0.0001299"
$ws.Range("E26").Value = "'-0.10%"
$ws.Range("D27").Value = "'0.0003396"
$ws.Range("E27").Value = "'-0.05%"
$ws.Range("D39").Value = "'0.01570"
$ws.Range("E39").Value = "'0.71%"
$ws.Range("D40").Value = "'0.04468"
$ws.Range("E40").Value = "'3.35%"
$ws.Range("D41").Value = "'0.007350"
$ws.Range("E41").Value = "'0.25%"
$ws.Range("D42").Value = "'0.009870"
$ws.Range("E42").Value = "'-1.24%"
$ws.Range("D43").Value = "'0.1330"
$ws.Range("E43").Value = "'2.15%"
$ws.Range("D44").Value = "'0.002249"
$ws.Range("E44").Value = "'11.71%"
$ws.Range("D45").Value = "'0.008769"
$ws.Range("E45").Value = "'-6.77%"
$ws.Range("D46").Value = "'0.00006113"
$ws.Range("E46").Value = "'3.71%"
$ws.Range("E47").Value = "'-0.05%"
$ws.Range("E49").Value = "'-31.04%"
$ws.Range("E50").Value = "'-0.05%"
$ws.Range("E51").Value = "'-0.05%"
